# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The underlying computation that produces the y_0_forecast (column C) and
# y_1_forecast (column E) values was corrected, which:
#   - slightly changed many of the previously-computed forecast values
#     (last-digit-level floating point differences from the recalculated
#     series), and
#   - removed the stray/invalid forecast value that had been written into
#     C2 (the first row has no prior-year y_0 to forecast from, so C2
#     should be blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear the erroneous C2 value, and update the recalculated E2 value.
$ws.Cells.Item(2, 3).Value = $null
$ws.Cells.Item(2, 5).Value = 5.080273296954396

# Row 3
$ws.Cells.Item(3, 3).Value = -3.942037578692492
$ws.Cells.Item(3, 5).Value = -1.648748515828491

# Row 4
$ws.Cells.Item(4, 3).Value = -2.839753013810498
$ws.Cells.Item(4, 5).Value = -1.632723506456923

# Row 5
$ws.Cells.Item(5, 3).Value = 4.960109259035428

# Row 7
$ws.Cells.Item(7, 3).Value = -2.700325749999499
$ws.Cells.Item(7, 5).Value = -0.3858735870725494

# Row 8
$ws.Cells.Item(8, 3).Value = 5.469647210234974
$ws.Cells.Item(8, 5).Value = 3.061326532789521

# Row 10
$ws.Cells.Item(10, 3).Value = 3.458696398997096
$ws.Cells.Item(10, 5).Value = 2.610227683091315

# Row 11
$ws.Cells.Item(11, 3).Value = 2.77241330895972

# Row 12
$ws.Cells.Item(12, 3).Value = 3.14581984265847

# Row 13
$ws.Cells.Item(13, 5).Value = 4.124307769579483

# Row 14
$ws.Cells.Item(14, 5).Value = 4.888255652935958

# Row 15
$ws.Cells.Item(15, 5).Value = 2.95288809451808

# Row 18
$ws.Cells.Item(18, 3).Value = -0.5744163079740128
$ws.Cells.Item(18, 5).Value = -0.6956477387308979

# Row 19
$ws.Cells.Item(19, 3).Value = -0.1892239049850142
